$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.027777443275147
$ws.Range("D2").Value = 1.034802810138583
$ws.Range("E2").Value = 1.036755587930699
$ws.Range("F2").Value = 1.045175645495356
$ws.Range("I2").Value = 1.034044921997951
$ws.Range("J2").Value = 1.032933430592864
$ws.Range("K2").Value = 1.037601237135002
$ws.Range("L2").Value = 1.039548421822541
$ws.Range("M2").Value = 1.047944622558245
$ws.Range("N2").Value = 1.034400314543426

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.02910187858249
$ws.Range("D3").Value = 1.035349979543521
$ws.Range("E3").Value = 1.037981007789417
$ws.Range("F3").Value = 1.046615171083977
$ws.Range("I3").Value = 1.034239994587896
$ws.Range("J3").Value = 1.03389627101164
$ws.Range("K3").Value = 1.037958367879821
$ws.Range("L3").Value = 1.040582406759339
$ws.Range("M3").Value = 1.04919389986713
$ws.Range("N3").Value = 1.035364522306036

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.029955198089136
$ws.Range("D4").Value = 1.035699477172608
$ws.Range("E4").Value = 1.038769278755952
$ws.Range("F4").Value = 1.04753768032756
$ws.Range("I4").Value = 1.034360504729651
$ws.Range("J4").Value = 1.034515204054786
$ws.Range("K4").Value = 1.038184202099592
$ws.Range("L4").Value = 1.041246260345409
$ws.Range("M4").Value = 1.04999281335286
$ws.Range("N4").Value = 1.035984334305097

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.03031306210619
$ws.Range("D5").Value = 1.035845315249325
$ws.Range("E5").Value = 1.039099561511138
$ws.Range("F5").Value = 1.047923370719239
$ws.Range("I5").Value = 1.034409799648978
$ws.Range("J5").Value = 1.034774431696239
$ws.Range("K5").Value = 1.03827788563578
$ws.Range("L5").Value = 1.041524106551226
$ws.Range("M5").Value = 1.050326424833028
$ws.Range("N5").Value = 1.036243930079534

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.030373098203734
$ws.Range("D6").Value = 1.035869738147932
$ws.Range("E6").Value = 1.039154952805458
$ws.Range("F6").Value = 1.047988005133397
$ws.Range("I6").Value = 1.034417996301208
$ws.Range("J6").Value = 1.034817900370803
$ws.Range("K6").Value = 1.038293541800386
$ws.Range("L6").Value = 1.041570685787011
$ws.Range("M6").Value = 1.050382307935886
$ws.Range("N6").Value = 1.036287460484602

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.029959983297485
$ws.Range("D7").Value = 1.035701430153224
$ws.Range("E7").Value = 1.038773696345654
$ws.Range("F7").Value = 1.047542842303623
$ws.Range("I7").Value = 1.034361168783341
$ws.Range("J7").Value = 1.034518671676022
$ws.Range("K7").Value = 1.038185458842597
$ws.Range("L7").Value = 1.041249977793014
$ws.Range("M7").Value = 1.049997279918407
$ws.Range("N7").Value = 1.035987806850752

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.028225809726419
$ws.Range("D8").Value = 1.034988672589767
$ws.Range("E8").Value = 1.037170694062133
$ws.Range("F8").Value = 1.045664002935994
$ws.Range("I8").Value = 1.034112031452158
$ws.Range("J8").Value = 1.033259678089199
$ws.Range("K8").Value = 1.037723019240659
$ws.Range("L8").Value = 1.039898944279025
$ws.Range("M8").Value = 1.048368785826935
$ws.Range("N8").Value = 1.034727025348628

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.025141370326956
$ws.Range("D9").Value = 1.033697763537451
$ws.Range("E9").Value = 1.034309916142613
$ws.Range("F9").Value = 1.042284024192095
$ws.Range("I9").Value = 1.033629206838212
$ws.Range("J9").Value = 1.03100951615009
$ws.Range("K9").Value = 1.036867869842036
$ws.Range("L9").Value = 1.03747801706583
$ws.Range("M9").Value = 1.045426238026363
$ws.Range("N9").Value = 1.032473667921467

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.023065218637142
$ws.Range("D10").Value = 1.032813591761838
$ws.Range("E10").Value = 1.032377863037827
$ws.Range("F10").Value = 1.039983268898548
$ws.Range("I10").Value = 1.033277779622235
$ws.Range("J10").Value = 1.029487628823372
$ws.Range("K10").Value = 1.036270613761782
$ws.Range("L10").Value = 1.035836453195587
$ws.Range("M10").Value = 1.043414695831009
$ws.Range("N10").Value = 1.030949619340183

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.022161364788882
$ws.Range("D11").Value = 1.032425121673025
$ws.Range("E11").Value = 1.031535219976805
$ws.Range("F11").Value = 1.038975550853281
$ws.Range("I11").Value = 1.033118573627305
$ws.Range("J11").Value = 1.028823351196992
$ws.Range("K11").Value = 1.036005528687767
$ws.Range("L11").Value = 1.035118954516422
$ws.Range("M11").Value = 1.042531658484288
$ws.Range("N11").Value = 1.030284398363379

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.021824888896377
$ws.Range("D12").Value = 1.032279979857933
$ws.Range("E12").Value = 1.03122130258515
$ws.Range("F12").Value = 1.038599495818703
$ws.Range("I12").Value = 1.033058377540404
$ws.Range("J12").Value = 1.028575803415059
$ws.Range("K12").Value = 1.035906089581178
$ws.Range("L12").Value = 1.034851426150483
$ws.Range("M12").Value = 1.042201834362181
$ws.Range("N12").Value = 1.030036499035202

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.021897098016509
$ws.Range("D13").Value = 1.032311151633158
$ws.Range("E13").Value = 1.031288680880304
$ws.Range("F13").Value = 1.038680240128699
$ws.Range("I13").Value = 1.03307133781701
$ws.Range("J13").Value = 1.028628939880647
$ws.Range("K13").Value = 1.035927463762702
$ws.Range("L13").Value = 1.034908858134378
$ws.Range("M13").Value = 1.042272665625452
$ws.Range("N13").Value = 1.030089710960665

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.02213356684867
$ws.Range("D14").Value = 1.032413141495155
$ws.Range("E14").Value = 1.031509290358337
$ws.Range("F14").Value = 1.038944501713727
$ws.Range("I14").Value = 1.033113619439144
$ws.Range("J14").Value = 1.0288029053152
$ws.Range("K14").Value = 1.03599732892117
$ws.Range("L14").Value = 1.03509686136554
$ws.Range("M14").Value = 1.042504432450147
$ws.Range("N14").Value = 1.03026392344609

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.022279164067372
$ws.Range("D15").Value = 1.032475868528908
$ws.Range("E15").Value = 1.031645092532291
$ws.Range("F15").Value = 1.039107090342795
$ws.Range("I15").Value = 1.033139530005039
$ws.Range("J15").Value = 1.028909984098719
$ws.Range("K15").Value = 1.036040245904281
$ws.Range("L15").Value = 1.035212561170535
$ws.Range("M15").Value = 1.04264698929476
$ws.Range("N15").Value = 1.030371154293764

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.023125100792991
$ws.Range("D16").Value = 1.032839254656841
$ws.Range("E16").Value = 1.032433657882081
$ws.Range("F16").Value = 1.04004990439288
$ws.Range("I16").Value = 1.033288197129189
$ws.Range("J16").Value = 1.029531602377741
$ws.Range("K16").Value = 1.036288070056085
$ws.Range("L16").Value = 1.035883929224102
$ws.Range("M16").Value = 1.043473045220638
$ws.Range("N16").Value = 1.030993655342042

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.023654422481999
$ws.Range("D17").Value = 1.033065691305821
$ws.Range("E17").Value = 1.032926674898548
$ws.Range("F17").Value = 1.040638220785842
$ws.Range("I17").Value = 1.033379566517304
$ws.Range("J17").Value = 1.029920103310821
$ws.Range("K17").Value = 1.036441789741885
$ws.Range("L17").Value = 1.036303260918313
$ws.Range("M17").Value = 1.04398797610819
$ws.Range("N17").Value = 1.031382707991009

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.023962698241253
$ws.Range("D18").Value = 1.033197226132506
$ws.Range("E18").Value = 1.033213660433471
$ws.Range("F18").Value = 1.040980269905518
$ws.Range("I18").Value = 1.033432182106289
$ws.Range("J18").Value = 1.030146199633057
$ws.Range("K18").Value = 1.036530827853114
$ws.Range("L18").Value = 1.036547205300447
$ws.Range("M18").Value = 1.044287167478048
$ws.Range("N18").Value = 1.031609125395956

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.024067733213478
$ws.Range("D19").Value = 1.033241984279748
$ws.Range("E19").Value = 1.033311416593932
$ws.Range("F19").Value = 1.041096712817952
$ws.Range("I19").Value = 1.033450007636076
$ws.Range("J19").Value = 1.030223206498527
$ws.Range("K19").Value = 1.036561081788341
$ws.Range("L19").Value = 1.036630275039352
$ws.Range("M19").Value = 1.044388987927526
$ws.Range("N19").Value = 1.031686241620008

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.023597679840849
$ws.Range("D20").Value = 1.033041452860111
$ws.Range("E20").Value = 1.032873839229191
$ws.Range("F20").Value = 1.040575214530187
$ws.Range("I20").Value = 1.033369833660306
$ws.Range("J20").Value = 1.029878473634195
$ws.Range("K20").Value = 1.036425361634874
$ws.Range("L20").Value = 1.036258337367457
$ws.Range("M20").Value = 1.043932848916261
$ws.Range("N20").Value = 1.031341019195468

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.022063953326726
$ws.Range("D21").Value = 1.032383131420041
$ws.Range("E21").Value = 1.031444351931683
$ws.Range("F21").Value = 1.038866731559072
$ws.Range("I21").Value = 1.033101197838014
$ws.Range("J21").Value = 1.028751699154055
$ws.Range("K21").Value = 1.035976782303958
$ws.Range("L21").Value = 1.035041527285853
$ws.Range("M21").Value = 1.042436233416507
$ws.Range("N21").Value = 1.030212644566323

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.021095325278542
$ws.Range("D22").Value = 1.031964318115221
$ws.Range("E22").Value = 1.030540233873554
$ws.Range("F22").Value = 1.037782441778778
$ws.Range("I22").Value = 1.032926161980319
$ws.Range("J22").Value = 1.028038585388999
$ws.Range("K22").Value = 1.035689101706459
$ws.Range("L22").Value = 1.034270577660712
$ws.Range("M22").Value = 1.041484685025931
$ws.Range("N22").Value = 1.029498518097925

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.021609226706962
$ws.Range("D23").Value = 1.032186804472012
$ws.Range("E23").Value = 1.031020034891776
$ws.Range("F23").Value = 1.038358208304548
$ws.Range("I23").Value = 1.033019534256294
$ws.Range("J23").Value = 1.028417066435238
$ws.Range("K23").Value = 1.035842142344975
$ws.Range("L23").Value = 1.034679835423222
$ws.Range("M23").Value = 1.041990126642446
$ws.Range("N23").Value = 1.029877536630662

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.023623320845873
$ws.Range("D24").Value = 1.033052406844958
$ws.Range("E24").Value = 1.032897715195081
$ws.Range("F24").Value = 1.040603687761877
$ws.Range("I24").Value = 1.033374233616508
$ws.Range("J24").Value = 1.029897285867994
$ws.Range("K24").Value = 1.036432786717582
$ws.Range("L24").Value = 1.036278638378357
$ws.Range("M24").Value = 1.043957762102707
$ws.Range("N24").Value = 1.031359858144796

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.025942224591453
$ws.Range("D25").Value = 1.034035641369704
$ws.Range("E25").Value = 1.03505383257525
$ws.Range("F25").Value = 1.043166122100426
$ws.Range("I25").Value = 1.033759228360486
$ws.Range("J25").Value = 1.031595037852785
$ws.Range("K25").Value = 1.037093725874251
$ws.Range("L25").Value = 1.038108708698885
$ws.Range("M25").Value = 1.046195680287956
$ws.Range("N25").Value = 1.033060021132141
